# The sheet's stat columns (runs/balls/fours/sixes) were re-ordered across
# the existing rows (a row-for-row shuffle of the C:F data, columns A/B -
# playerName/teamName - stay put). Each target cell already held a
# numeric-looking value stored as TEXT (OOXML t="str"), so values are
# (re)written with a leading apostrophe to keep them text instead of
# letting Excel auto-convert them to real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'10"

$ws.Range("C3").Value  = "'2"
$ws.Range("D3").Value  = "'3"
$ws.Range("E3").Value  = "'0"
$ws.Range("F3").Value  = "'0"

$ws.Range("C4").Value  = "'16"
$ws.Range("D4").Value  = "'7"
$ws.Range("E4").Value  = "'1"
$ws.Range("F4").Value  = "'1"

$ws.Range("C5").Value  = "'1"
$ws.Range("D5").Value  = "'4"

$ws.Range("C6").Value  = "'13"
$ws.Range("D6").Value  = "'3"
$ws.Range("F6").Value  = "'2"

$ws.Range("C7").Value  = "'27"
$ws.Range("D7").Value  = "'8"
$ws.Range("E7").Value  = "'0"
$ws.Range("F7").Value  = "'4"

$ws.Range("C8").Value  = "'6"
$ws.Range("D8").Value  = "'9"

$ws.Range("C9").Value  = "'2"
$ws.Range("D9").Value  = "'4"
$ws.Range("F9").Value  = "'0"

$ws.Range("C10").Value = "'24"
$ws.Range("D10").Value = "'11"
$ws.Range("E10").Value = "'3"

$ws.Range("C11").Value = "'6"
$ws.Range("F11").Value = "'1"
